$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 140; this shifts the existing rows 140-191
# down to 141-192 (dimension grows from A1:R191 to A1:R192) and is the
# mechanism behind every other cell delta shown in the diff.
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with its data.
$ws.Cells.Item(140, 1).Value = 6
$ws.Cells.Item(140, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(140, 3).Value = "Metropolitana"
$ws.Cells.Item(140, 4).Value = 44559
$ws.Cells.Item(140, 5).Value = 13
$ws.Cells.Item(140, 6).Value = 100112022
$ws.Cells.Item(140, 7).Value = "Arveja Verde"
$ws.Cells.Item(140, 8).Value = "Perfection"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 180
$ws.Cells.Item(140, 11).Value = 25000
$ws.Cells.Item(140, 12).Value = 27000
$ws.Cells.Item(140, 13).Value = 25889
$ws.Cells.Item(140, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(140, 15).Value = "Carahue"
$ws.Cells.Item(140, 16).Value = 1036
$ws.Cells.Item(140, 17).Value = 25
$ws.Cells.Item(140, 18).Value = "Hortaliza"
